$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(
    149.0429963588715,72.59762640476231,22.89976842880253,13.34329990863802,-15.10997164770959,-24.45587922990319,-45.75726741790768,-47.4942824053764,-18.41118627578018,-35.23416106343268,
    -25.05642833709715,-16.47099995128809,3.641371173858658,-56.43801583290099,-15.20868131726979,-34.42504441499707,-21.4366173329949,-7.029076403379406,-33.01469397902486,-0.9288858437537897,
    -5.880547896027565,-4.745866158604599,-30.71236153125761,-33.23559818506237,-0.01870143532750035,-8.384757556319187,6.341149383783375,-3.99053210318084,-8.64949123919007,24.11718803167344,
    97.40048725128179,86.501018652916,179.4538644313812,82.37207330226903,24.09675856113438,-9.33716943860054,-4.505535907149294,-49.52012155771251,-50.27269371032713,-96.71495047092435,
    -13.99249571204183,-28.1770614224672,-7.583838030099855,-42.22100687146184,-32.05286960482596,-51.21834075212476,-14.81343137547373,-11.29744734346866,-30.22298480749129,3.852774301767397,
    0.567045495510115,-20.4533401933312,12.5613668084145,20.35799324750906,7.602782613039039,2.1755326533318,31.6344085931778,-13.64262488484383,11.34847420930868,13.66719096660614,
    19.02034145832067,-0.1074714422225611,-22.0848701149225,19.49174594402314,-2.190245863199209,12.53134353399281,-4.57257579922674,-23.80666690200565,27.6320040917397,24.12849726200108,
    3.159887729883224,-3.742797865867601,-26.7398428714275,-24.07472753524777,-23.52222527086732,-38.74108715534209,-13.71795959234237,-43.6070946204662,-14.32488690674302,-14.9602619856596,
    -2.814172623157447,-10.66125015139579,-11.13046357721089,11.77696985244756,-12.40066524565219,-24.49549548029896,18.84307299852372,-13.39605634272095,-20.33017348527906,-15.45381033346052,
    -0.7565550303459077,-23.89186241388319,-15.75048783153295,-59.73075371265409,-42.94979486823078,-4.528704926967578,4.412129164934186,-22.3001544779539,-8.920051742196051,-28.23429225444792,
    -15.23089099809522,19.08066287279132,-41.1773470520973,-35.91318844676016,-10.09366918921467,0.489859584569956,-45.51545852184293,30.553262732029,18.23998469591146,-9.241818386912314,
    -3.182894927263249,-2.99157685995101,3.886980170011554,4.373742358684581,21.03634836912158,-10.51229450702664,-26.18473600685593,17.70338738203054,-31.53196446895598,-30.83785838782785,
    -37.05636949300762,-27.41996723353861,-38.02799464821814,-24.30399444520469,0.6155774927139532,-37.97131271958349,-11.45558863282201,1.53337203621868,-14.32913697510958,7.688125078678183,
    73.63517460346225,24.76675180673601,-6.451350250244104,3.133975849151625,-11.44016511112449,-15.42515724763274,7.339352058172267,30.15458504199984,10.4556974053383,-11.11976989686485,
    -21.52064756363629,4.011601927280481,2.668739281892783,-5.870951608419375,-14.67996755525468,9.167331010103226,74.44411407947541
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Add new rows 105-148 for column A (index values), copying the style from an existing A-column cell (e.g. A2)
$styleSource = $ws.Range("A2")
for ($row = 105; $row -le 148; $row++) {
    $idx = $row - 2
    $cellA = $ws.Cells.Item($row, 1)
    $styleSource.Copy($cellA)
    $cellA.Value = $idx
}

$ws.Range("A1").Select() | Out-Null